# Edit description (from the commit's OOXML diff):
#   1. Three tables (on slides 14, 15 and 16) have their table style
#      switched from the custom "Table_0" style
#      ({F573D811-D1D7-416D-A4BE-BE076CB7CDD4}, defined in
#      ppt/tableStyles.xml) to the built-in style
#      {42128279-9664-4AFC-A7CB-5FE5C3B4959F}.
#   2. The presentation's applied design is switched from the custom
#      "Integral" / "Red Violet" theme back to the standard
#      "Office Theme" / "Office" colour palette (the font scheme and
#      format scheme are identical between the two themes, so only the
#      12 theme colours actually change).

$p = $ppt.ActivePresentation

# --- 1. Table styles -------------------------------------------------
$newTableStyle = "{42128279-9664-4AFC-A7CB-5FE5C3B4959F}"

foreach ($slideIndex in 14..16) {
    $slide = $p.Slides.Item($slideIndex)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($newTableStyle)
        }
    }
}

# --- 2. Theme colours --------------------------------------------------
# The presentation theme (shared by every slide/layout/master) exposes
# its 12 colour slots through ThemeColorScheme, indexed in the standard
# order: Dark1, Light1, Dark2, Light2, Accent1-6, Hyperlink, Followed
# Hyperlink. Restore the stock "Office" palette values.
$officeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$themeColors = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $officeColors.Count; $i++) {
    $themeColors.Item($i).RGB = $officeColors[$i - 1]
}
